$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("BonusPower", 2,                 960,  106, 23.32, 46.64, "win"),
    @("BonusPower", 2,                 379,  225, 20,    40,    "win"),
    @("BonusPower", 2.369999999999992, 0,    396, 20,    0,     "lose"),
    @("BonusPower", 2.369999999999992, 0,    397, 20,    0,     "lose"),
    @("SkipBoss",   2,                 746,  172, 40,    80,    "win"),
    @("BonusPower", 2,                 1030, 113, 100,   200,   "win"),
    @("BonusPower", 0.03,              -170, 3,   30,    0,     "win"),
    @("SkipBoss",   2,                 1047, 120, 20,    40,    "win")
)

$startRow = 13
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
}
